$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 309, pushing the existing rows 309:383 down to 311:385.
$ws.Rows.Item(309).Resize(2, 1).EntireRow.Insert()

# New row 309 ("Primera")
$ws.Cells.Item(309, 1).Value = 8
$ws.Cells.Item(309, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = 44508
$ws.Cells.Item(309, 5).Value = 4
$ws.Cells.Item(309, 6).Value = 100112043
$ws.Cells.Item(309, 7).Value = "Pepino ensalada"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 800
$ws.Cells.Item(309, 11).Value = 7000
$ws.Cells.Item(309, 12).Value = 7500
$ws.Cells.Item(309, 13).Value = 7250
$ws.Cells.Item(309, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(309, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(309, 16).Value = 121
$ws.Cells.Item(309, 17).Value = 60
$ws.Cells.Item(309, 18).Value = "Hortaliza"

# New row 310 ("Segunda")
$ws.Cells.Item(310, 1).Value = 8
$ws.Cells.Item(310, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(310, 3).Value = "Coquimbo"
$ws.Cells.Item(310, 4).Value = 44508
$ws.Cells.Item(310, 5).Value = 4
$ws.Cells.Item(310, 6).Value = 100112043
$ws.Cells.Item(310, 7).Value = "Pepino ensalada"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Segunda"
$ws.Cells.Item(310, 10).Value = 400
$ws.Cells.Item(310, 11).Value = 4500
$ws.Cells.Item(310, 12).Value = 5000
$ws.Cells.Item(310, 13).Value = 4750
$ws.Cells.Item(310, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(310, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(310, 16).Value = 48
$ws.Cells.Item(310, 17).Value = 100
$ws.Cells.Item(310, 18).Value = "Hortaliza"
